$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-5 changed, row 6 added)
$data = @(
    @(1, 7, 4, 6, 9, -1, 5, 21, 5),
    @(2, 5, 4, 3, 8, -2, 4, 32, 5),
    @(3, 9, 0, 4, 1, -5, 1, 65, 5),
    @(4, 8, 3, 5, 6, -3, 3, 43, 5),
    @(5, 6, 1, 2, 3, -4, 2, 54, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select() | Out-Null
